$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.623.21'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '2.279.23'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '2.296.95'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0968'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.54%  '
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.340'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.42%  '
$ws.Range("E13").Value = '  +3.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.95%  '
$ws.Range("D15").Value = '2.686.36'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").Value = '54.670.96'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '2.292.96'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '307.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.995'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.08'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.20%  '
$ws.Range("D30").Value = '0.0₃0703'
$ws.Range("E30").Value = '  +2.78%  '
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.90%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.995'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.909'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.375'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("E41").Value = '  +1.07%  '
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '128.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.99%  '
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '250.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.46%  '
$ws.Range("E46").Value = '  +1.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0906'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.550'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("E51").Value = '  +0.39%  '
